$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '26.949.44'
$ws.Range('E2').Value = '  -2.76%  '
Set-TextValue 'D3' '1.860.11'
$ws.Range('E3').Value = '  -2.26%  '
Set-TextValue 'D4' '0.9995'
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue 'D5' '305.75'
$ws.Range('E5').Value = '  -2.03%  '
Set-TextValue 'D6' '0.9993'
Set-TextValue 'D7' '0.5035'
$ws.Range('E7').Value = '  -2.91%  '
$ws.Range('E8').Value = '  -1.87%  '
Set-TextValue 'D9' '0.07121'
$ws.Range('E9').Value = '  -1.57%  '
Set-TextValue 'D10' '0.8846'
$ws.Range('E10').Value = '  -1.24%  '
Set-TextValue 'D11' '20.54'
$ws.Range('E11').Value = '  -2.73%  '
Set-TextValue 'D12' '1.881.07'
$ws.Range('E12').Value = '  -0.93%  '
Set-TextValue 'D13' '0.07566'
$ws.Range('E13').Value = '  -1.02%  '
Set-TextValue 'D14' '5.282'
$ws.Range('E14').Value = '  -2.95%  '
Set-TextValue 'D15' '89.02'
$ws.Range('E15').Value = '  -3.35%  '
Set-TextValue 'D16' '1.000'
$ws.Range('E16').Value = '  -0.01%  '
Set-TextValue 'D17' '0.000008353'
$ws.Range('E17').Value = '  -4.14%  '
Set-TextValue 'D18' '14.06'
$ws.Range('E18').Value = '  -2.86%  '
Set-TextValue 'D19' '0.9991'
Set-TextValue 'D20' '26.987.75'
$ws.Range('E20').Value = '  -2.70%  '
Set-TextValue 'D21' '5.022'
$ws.Range('E21').Value = '  -2.16%  '
Set-TextValue 'D22' '2.119.53'
$ws.Range('E22').Value = '  -1.85%  '
$ws.Range('E23').Value = '  -3.51%  '
Set-TextValue 'D24' '6.453'
$ws.Range('E24').Value = '  -1.89%  '
$ws.Range('E25').Value = '  -0.81%  '
Set-TextValue 'D26' '146.91'
$ws.Range('E26').Value = '  -4.41%  '
Set-TextValue 'D27' '17.92'
$ws.Range('E27').Value = '  -2.12%  '
$ws.Range('E28').Value = '  -4.23%  '
Set-TextValue 'D29' '112.39'
$ws.Range('E29').Value = '  -2.06%  '
Set-TextValue 'D30' '4.641'
$ws.Range('E30').Value = '  -4.24%  '
Set-TextValue 'D31' '4.657'
$ws.Range('E31').Value = '  -3.34%  '
Set-TextValue 'D32' '0.09033'
$ws.Range('E32').Value = '  +0.68%  '
Set-TextValue 'D33' '0.05107'
$ws.Range('E33').Value = '  -3.17%  '
Set-TextValue 'D34' '3.041'
$ws.Range('E34').Value = '  -4.21%  '
Set-TextValue 'D35' '1.147'
$ws.Range('E35').Value = '  -7.48%  '
Set-TextValue 'D36' '0.7216'
$ws.Range('E36').Value = '  -7.27%  '
Set-TextValue 'D37' '0.02029'
$ws.Range('E37').Value = '  -2.64%  '
Set-TextValue 'D38' '3.033'
$ws.Range('E38').Value = '  -0.66%  '
Set-TextValue 'D39' '2.458'
$ws.Range('E39').Value = '  -6.13%  '
Set-TextValue 'D40' '1.072'
$ws.Range('E40').Value = '  -1.64%  '
Set-TextValue 'D41' '0.5267'
$ws.Range('E41').Value = '  -4.18%  '
Set-TextValue 'D42' '6.514'
$ws.Range('E42').Value = '  -2.10%  '
Set-TextValue 'D43' '115.08'
$ws.Range('E43').Value = '  +1.76%  '
Set-TextValue 'D44' '8.233'
$ws.Range('E44').Value = '  -3.09%  '
Set-TextValue 'D45' '0.1459'
$ws.Range('E45').Value = '  -3.00%  '
Set-TextValue 'D46' '0.9991'
Set-TextValue 'D47' '0.4589'
$ws.Range('E47').Value = '  -4.30%  '
Set-TextValue 'D48' '9.954'
$ws.Range('E48').Value = '  -4.77%  '
Set-TextValue 'D49' '1.558'
$ws.Range('E49').Value = '  -3.59%  '
Set-TextValue 'D50' '36.43'
$ws.Range('E50').Value = '  -0.94%  '
Set-TextValue 'D51' '63.86'
$ws.Range('E51').Value = '  -4.05%  '
